$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.9
$ws.Cells.Item(2, 9).Value = 2.6
$ws.Cells.Item(2, 10).Value = 3.75
$ws.Cells.Item(2, 11).Value = 1.91
$ws.Cells.Item(2, 12).Value = 3.4
$ws.Cells.Item(2, 15).Value = 1.5
$ws.Cells.Item(2, 16).Value = 2.5
$ws.Cells.Item(2, 17).Value = 1.93
$ws.Cells.Item(2, 18).Value = 1.93
$ws.Cells.Item(2, 20).Value = 1.5
$ws.Cells.Item(2, 25).Value = 1.57
$ws.Cells.Item(2, 26).Value = 2.25
$ws.Cells.Item(2, 30).Value = 13
$ws.Cells.Item(2, 33).Value = 29
$ws.Cells.Item(2, 37).Value = 19
$ws.Cells.Item(2, 39).Value = 501
$ws.Cells.Item(2, 40).Value = 6.5
$ws.Cells.Item(2, 41).Value = 11
$ws.Cells.Item(2, 43).Value = 26
$ws.Cells.Item(2, 44).Value = 23

# Row 3
$ws.Cells.Item(3, 7).Value = 1.62
$ws.Cells.Item(3, 8).Value = 3.8
$ws.Cells.Item(3, 9).Value = 5.5
$ws.Cells.Item(3, 11).Value = 2.2
$ws.Cells.Item(3, 12).Value = 5.5
$ws.Cells.Item(3, 17).Value = 1.49
$ws.Cells.Item(3, 18).Value = 2.65
$ws.Cells.Item(3, 21).Value = 2.7
$ws.Cells.Item(3, 22).Value = 1.47
$ws.Cells.Item(3, 30).Value = 7.5
$ws.Cells.Item(3, 32).Value = 12
$ws.Cells.Item(3, 34).Value = 26
$ws.Cells.Item(3, 42).Value = 17
$ws.Cells.Item(3, 43).Value = 51
$ws.Cells.Item(3, 45).Value = 41

# Row 4
$ws.Cells.Item(4, 7).Value = 2.1
$ws.Cells.Item(4, 8).Value = 3.4
$ws.Cells.Item(4, 11).Value = 1.91
$ws.Cells.Item(4, 13).Value = 1.1
$ws.Cells.Item(4, 14).Value = 7
$ws.Cells.Item(4, 15).Value = 1.53
$ws.Cells.Item(4, 16).Value = 2.38
$ws.Cells.Item(4, 17).Value = 2.03
$ws.Cells.Item(4, 18).Value = 1.83
$ws.Cells.Item(4, 19).Value = 2.6
$ws.Cells.Item(4, 20).Value = 1.48
$ws.Cells.Item(4, 21).Value = 4.3
$ws.Cells.Item(4, 22).Value = 1.22
$ws.Cells.Item(4, 23).Value = 5.5
$ws.Cells.Item(4, 24).Value = 1.14
$ws.Cells.Item(4, 31).Value = 10
$ws.Cells.Item(4, 32).Value = 19
$ws.Cells.Item(4, 35).Value = 6.5
$ws.Cells.Item(4, 40).Value = 7.5
$ws.Cells.Item(4, 44).Value = 41
$ws.Cells.Item(4, 45).Value = 51

# Row 9
$ws.Cells.Item(9, 7).Value = 3.1
$ws.Cells.Item(9, 8).Value = 3.1
$ws.Cells.Item(9, 9).Value = 2.25
$ws.Cells.Item(9, 10).Value = 3.55
$ws.Cells.Item(9, 11).Value = 2.1
$ws.Cells.Item(9, 12).Value = 2.77
$ws.Cells.Item(9, 19).Value = 1.85
$ws.Cells.Item(9, 31).Value = 10.75
$ws.Cells.Item(9, 33).Value = 27
$ws.Cells.Item(9, 36).Value = 6.1
$ws.Cells.Item(9, 37).Value = 13
$ws.Cells.Item(9, 39).Value = 450
$ws.Cells.Item(9, 40).Value = 7.9
$ws.Cells.Item(9, 41).Value = 11.25
$ws.Cells.Item(9, 43).Value = 23
$ws.Cells.Item(9, 44).Value = 18

# Row 10
$ws.Cells.Item(10, 7).Value = 1.85
$ws.Cells.Item(10, 8).Value = 3.3
$ws.Cells.Item(10, 9).Value = 3.9
$ws.Cells.Item(10, 10).Value = 2.63
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 5
$ws.Cells.Item(10, 14).Value = 7.5
$ws.Cells.Item(10, 15).Value = 1.44
$ws.Cells.Item(10, 16).Value = 2.63
$ws.Cells.Item(10, 17).Value = 1.78
$ws.Cells.Item(10, 18).Value = 2.03
$ws.Cells.Item(10, 19).Value = 2.35
$ws.Cells.Item(10, 20).Value = 1.57
$ws.Cells.Item(10, 23).Value = 4.5
$ws.Cells.Item(10, 24).Value = 1.18
$ws.Cells.Item(10, 25).Value = 1.53
$ws.Cells.Item(10, 26).Value = 2.38
$ws.Cells.Item(10, 30).Value = 8
$ws.Cells.Item(10, 31).Value = 9.5
$ws.Cells.Item(10, 32).Value = 15
$ws.Cells.Item(10, 33).Value = 19
$ws.Cells.Item(10, 35).Value = 7.5
$ws.Cells.Item(10, 36).Value = 6.5
$ws.Cells.Item(10, 40).Value = 9
$ws.Cells.Item(10, 41).Value = 19
$ws.Cells.Item(10, 42).Value = 15
$ws.Cells.Item(10, 43).Value = 41
$ws.Cells.Item(10, 45).Value = 41

# Row 11
$ws.Cells.Item(11, 7).Value = 2.2
$ws.Cells.Item(11, 8).Value = 2.87
$ws.Cells.Item(11, 9).Value = 3.25
$ws.Cells.Item(11, 13).Value = 1.1
$ws.Cells.Item(11, 14).Value = 7
$ws.Cells.Item(11, 19).Value = 2.6
$ws.Cells.Item(11, 20).Value = 1.48
$ws.Cells.Item(11, 27).Value = 2.1
$ws.Cells.Item(11, 28).Value = 1.67
$ws.Cells.Item(11, 33).Value = 23

# Row 16
$ws.Cells.Item(16, 7).Value = 11
$ws.Cells.Item(16, 9).Value = 1.17
$ws.Cells.Item(16, 10).Value = 10
$ws.Cells.Item(16, 14).Value = 15
$ws.Cells.Item(16, 19).Value = 1.33
$ws.Cells.Item(16, 20).Value = 3.25
$ws.Cells.Item(16, 27).Value = 2
$ws.Cells.Item(16, 28).Value = 1.73
$ws.Cells.Item(16, 32).Value = 151
$ws.Cells.Item(16, 33).Value = 67
$ws.Cells.Item(16, 35).Value = 26
$ws.Cells.Item(16, 37).Value = 29
$ws.Cells.Item(16, 38).Value = 67
$ws.Cells.Item(16, 39).Value = 700
$ws.Cells.Item(16, 40).Value = 11
$ws.Cells.Item(16, 41).Value = 7.5
$ws.Cells.Item(16, 43).Value = 7.5
$ws.Cells.Item(16, 45).Value = 26
